# B6-PowerPoint.pptx edit script
#
# 1) Re-style the three tables (slides 14, 15, 16) from the custom
#    "Table_0" style {1FBFB0C8-2BAB-49FB-AC28-D1839D21F73A} to the
#    built-in table style {0A228784-0033-41D0-8C04-6CA0306D8CA8}.
#
# 2) Swap the presentation's colour theme over to the stock "Office"
#    palette (previously the deck used the "Integral"/"Red Violet"
#    palette). This is done by rewriting the 12 theme colour slots
#    (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink) through
#    Slide.ThemeColorScheme, which is the scriptable surface for the
#    presentation's main theme (ppt/theme/theme1.xml).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# 1) Table styles
# ---------------------------------------------------------------
$newTableStyle = "{0A228784-0033-41D0-8C04-6CA0306D8CA8}"
$tableSlideNumbers = @(14, 15, 16)

foreach ($slideNum in $tableSlideNumbers) {
    $slide = $p.Slides.Item($slideNum)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}

# ---------------------------------------------------------------
# 2) Theme colours -> Office palette
# ---------------------------------------------------------------
# Slot order matches DrawingML's <a:clrScheme>: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink. COM hands back/accepts colours as
# OLE COLORREF (0xBBGGRR) integers, hence the byte-swapped values
# below relative to the RRGGBB hex they represent.
$officeColorsBGR = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeColorsBGR[$i - 1]
}
